$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text interpretation
# (so numeric-looking strings like "213.62" are not silently converted
# to numbers) and then clear the formatting so the cell's style stays
# the same as before the edit (no stray NumberFormat="@" on the cell).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Rows with only Price (D) and/or Volume(1h) (E) updates
$rows = @(
    @{Row=2;  D="27.900.78";  E="  +1.23%  "},
    @{Row=3;  D="1.643.30";   E="  +1.21%  "},
    @{Row=5;  D="213.62";     E="  +0.97%  "},
    @{Row=6;  E="  -0.29%  "},
    @{Row=7;  E="  +0.06%  "},
    @{Row=8;  D="23.56";      E="  +1.31%  "},
    @{Row=9;  E="  +0.14%  "},
    @{Row=10; D="0.0617";     E="  +0.82%  "},
    @{Row=11; D="0.0878";     E="  -1.31%  "},
    @{Row=12; D="1.876.07";   E="  +1.24%  "},
    @{Row=13; D="1.646.11";   E="  +1.36%  "},
    @{Row=14; E="  +4.11%  "},
    @{Row=15; E="  +0.46%  "},
    @{Row=16; D="65.86";      E="  +1.00%  "},
    @{Row=17; D="27.885.90";  E="  +1.29%  "},
    @{Row=18; D="230.70"},
    @{Row=19; E="  +0.75%  "},
    @{Row=20; D="7.62";       E="  +0.69%  "},
    @{Row=21; E="  +0.06%  "},
    @{Row=22; E="  +4.77%  "},
    @{Row=23; E="  +1.59%  "},
    @{Row=24; E="  +2.52%  "},
    @{Row=25; D="152.68";     E="  +1.67%  "},
    @{Row=26; D="6.93";       E="  +0.73%  "},
    @{Row=27; E="  +0.88%  "},
    @{Row=28; D="15.73";      E="  +1.03%  "},
    @{Row=29; E="  +0.03%  "},
    @{Row=30; E="  +1.19%  "},
    @{Row=31; E="  +0.31%  "},
    @{Row=32; E="  +1.87%  "},
    @{Row=33; D="1.434.44";   E="  -2.81%  "},
    @{Row=34; D="3.09";       E="  +0.39%  "},
    @{Row=35; E="  +1.82%  "},
    @{Row=36; E="  +0.19%  "},
    @{Row=37; D="0.886";      E="  +1.64%  "},
    @{Row=40; D="0.559";      E="  +0.77%  "},
    @{Row=41; E="  +1.72%  "},
    @{Row=44; E="  +0.42%  "},
    @{Row=45; E="  +3.13%  "},
    @{Row=46; D="1.80";       E="  +3.22%  "},
    @{Row=47; E="  +0.07%  "},
    @{Row=48; D="1.784.88";   E="  +1.16%  "},
    @{Row=49; D="89.13";      E="  +2.04%  "},
    @{Row=50; E="  +0.05%  "},
    @{Row=51; E="  +0.53%  "}
)

foreach ($r in $rows) {
    if ($r.ContainsKey("D")) {
        Set-TextValue ("D" + $r.Row) $r.D
    }
    if ($r.ContainsKey("E")) {
        Set-TextValue ("E" + $r.Row) $r.E
    }
}

# Rows 38/39 effectively swap VeChain and TrustWalletToken (with updated
# price/volume figures for each).
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "0.930"
Set-TextValue "E38" "  -2.54%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0168"
Set-TextValue "E39" "  +0.66%  "

# Rows 42/43 effectively swap PaxDollar and Aave (with updated
# price/volume figures for each).
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D42" "68.68"
Set-TextValue "E42" "  +1.41%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D43" "1.00"
Set-TextValue "E43" "  -0.01%  "
